# Edit: Finitos2, Taller 2 2022b
# Applies cell-value and view/selection changes per the target diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("xnod")
$ws2 = $wb.Worksheets.Item("LaG_mat")
$ws3 = $wb.Worksheets.Item("restric")
$ws4 = $wb.Worksheets.Item("carga_distr")
$ws5 = $wb.Worksheets.Item("carga_punt")
$ws6 = $wb.Worksheets.Item("prop_mat")
$ws7 = $wb.Worksheets.Item("config")
$ws8 = $wb.Worksheets.Item("varios")

# --- carga_distr (sheet4): populate column F (EF type per bar: ER_PH/EE/ER_HP cycling)
# and set D/E (b1,b2) offsets to -20 for the rigid-link bars.
$ws4.Range("F2").Value = "ER_PH"
$ws4.Range("F3").Value = "EE"
$ws4.Range("F4").Value = "ER_HP"
$ws4.Range("F5").Value = "ER_PH"
$ws4.Range("F6").Value = "EE"
$ws4.Range("F7").Value = "ER_HP"
$ws4.Range("F8").Value = "ER_PH"
$ws4.Range("F9").Value = "EE"
$ws4.Range("F10").Value = "ER_HP"
$ws4.Range("D11").Value = -20
$ws4.Range("E11").Value = -20
$ws4.Range("F11").Value = "ER_PH"
$ws4.Range("F12").Value = "EE"
$ws4.Range("D13").Value = -20
$ws4.Range("E13").Value = -20
$ws4.Range("F13").Value = "ER_HP"
$ws4.Range("D14").Value = -20
$ws4.Range("E14").Value = -20
$ws4.Range("F14").Value = "ER_PH"
$ws4.Range("F15").Value = "EE"
$ws4.Range("D16").Value = -20
$ws4.Range("E16").Value = -20
$ws4.Range("F16").Value = "ER_HP"
$ws4.Range("F17").Value = "ER_PH"
$ws4.Range("F18").Value = "EE"
$ws4.Range("F19").Value = "ER_HP"
$ws4.Range("F20").Value = "ER_PH"
$ws4.Range("F21").Value = "EE"
$ws4.Range("F22").Value = "ER_HP"
$ws4.Range("F23").Value = "ER_PH"
$ws4.Range("F24").Value = "EE"
$ws4.Range("F25").Value = "ER_HP"
$ws4.Range("D26").Value = -20
$ws4.Range("E26").Value = -20
$ws4.Range("F26").Value = "ER_PH"
$ws4.Range("F27").Value = "EE"
$ws4.Range("D28").Value = -20
$ws4.Range("E28").Value = -20
$ws4.Range("F28").Value = "ER_HP"
$ws4.Range("D29").Value = -20
$ws4.Range("E29").Value = -20
$ws4.Range("F29").Value = "ER_PH"
$ws4.Range("F30").Value = "EE"
$ws4.Range("D31").Value = -20
$ws4.Range("E31").Value = -20
$ws4.Range("F31").Value = "ER_HP"
$ws4.Range("F32").Value = "ER_PH"
$ws4.Range("F33").Value = "EE"
$ws4.Range("F34").Value = "ER_HP"
$ws4.Range("F35").Value = "ER_PH"
$ws4.Range("F36").Value = "EE"
$ws4.Range("F37").Value = "ER_HP"
$ws4.Range("F38").Value = "ER_PH"
$ws4.Range("F39").Value = "EE"
$ws4.Range("F40").Value = "ER_HP"
$ws4.Range("D41").Value = -20
$ws4.Range("E41").Value = -20
$ws4.Range("F41").Value = "ER_PH"
$ws4.Range("F42").Value = "EE"
$ws4.Range("D43").Value = -20
$ws4.Range("E43").Value = -20
$ws4.Range("F43").Value = "ER_HP"
$ws4.Range("D44").Value = -20
$ws4.Range("E44").Value = -20
$ws4.Range("F44").Value = "ER_PH"
$ws4.Range("F45").Value = "EE"
$ws4.Range("D46").Value = -20
$ws4.Range("E46").Value = -20
$ws4.Range("F46").Value = "ER_HP"
$ws4.Range("F47").Value = "ER_PH"
$ws4.Range("F48").Value = "EE"
$ws4.Range("F49").Value = "ER_HP"
$ws4.Range("F50").Value = "ER_PH"
$ws4.Range("F51").Value = "EE"
$ws4.Range("F52").Value = "ER_HP"
$ws4.Range("F53").Value = "ER_PH"
$ws4.Range("F54").Value = "EE"
$ws4.Range("F55").Value = "ER_HP"
$ws4.Range("D56").Value = -20
$ws4.Range("E56").Value = -20
$ws4.Range("F56").Value = "ER_PH"
$ws4.Range("F57").Value = "EE"
$ws4.Range("D58").Value = -20
$ws4.Range("E58").Value = -20
$ws4.Range("F58").Value = "ER_HP"
$ws4.Range("D59").Value = -20
$ws4.Range("E59").Value = -20
$ws4.Range("F59").Value = "ER_PH"
$ws4.Range("F60").Value = "EE"
$ws4.Range("D61").Value = -20
$ws4.Range("E61").Value = -20
$ws4.Range("F61").Value = "ER_HP"

# carga_punt (sheet5): clear sample load value C2 (10 -> 0)
$ws5.Range("C2").Value = 0

# config (sheet7): esc_faxial 0.05 -> 0.01
$ws7.Range("B3").Value = 0.01

# --- Update each sheet's view selection (activeCell/sqref) to match the saved state,
# finishing on "config" so it ends up as the active tab (matches activeTab="6").

$ws2.Range("E61").Select()
$ws4.Range("H14").Select()
$ws5.Range("C3").Select()
$ws6.Range("G23").Select()
$ws8.Range("E30").Select()
$ws7.Range("B3").Select()
